$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert columns to make room for the new fields ---
# New column A: "Codigo" (before the old Pedido column)
$ws.Columns.Item(1).Insert()
# New column D: "Material" (after old Item column, which is now C)
$ws.Columns.Item(4).Insert()
# New columns F, G, H: "Data de Remessa", "Fornecedor", "Follow-up"
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(7).Insert()
$ws.Columns.Item(8).Insert()

# Match column widths for the newly inserted columns to their neighbours /
# to the width needed to fit their contents (mirrors Excel's own best-fit sizing)
$ws.Range("D1").ColumnWidth = $ws.Range("C1").ColumnWidth
$ws.Range("F1").ColumnWidth = 15.5
$ws.Range("G1").ColumnWidth = $ws.Range("C1").ColumnWidth
$ws.Range("H1").ColumnWidth = 22.166666666666668

# --- Material (column D) ---
$ws.Range("D1").Value = "Material"
$ws.Range("D2").Value = "Sabonete"
$ws.Range("D3").Value = "Detergente"
$ws.Range("D5").Value = "Acido tipo 2"
$ws.Range("D6").Value = "Etanol"
$ws.Range("D4").Value = "Acido TIpo 1"

# --- Codigo (column A): header + formulas ---
$ws.Range("A1").Value = "Codigo"
$ws.Range("A2").Formula = "=_xlfn.CONCAT(RIGHT(B2,5),C2)"
$ws.Range("A3:A6").Formula = "=_xlfn.CONCAT(RIGHT(B3,5),C3)"

# --- Data de Remessa / Fornecedor / Follow-up headers ---
$ws.Range("F1").Value = "Data de Remessa"
$ws.Range("G1").Value = "Fornecedor"
$ws.Range("H1").Value = "Follow-up"

# --- Fornecedor (column G) ---
$ws.Range("G2:G6").Value = "Merck"

# --- Follow-up (column H) ---
$ws.Range("H2").Value = "Aguardando importação"
$ws.Range("H4").Value = "Fora de estoque"

# --- Data de Remessa values (column F), with short-date format ---
$ws.Range("F2").Value = 45667
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Copy()
$ws.Range("F3:F6").PasteSpecial(-4122)
$ws.Range("F3").Value = 45677
$ws.Range("F4").Value = 45669
$ws.Range("F5").Value = 45677
$ws.Range("F6").Value = 45677

# --- Selection matches the committed workbook state ---
$ws.Range("F1").Select()
